$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings are created in this order as they are first referenced:
# "44 u=4 6u6nPBC", "44 u=4 7u7nPBC", "44 u=4 8u8nPBC", "ED", "t'", "t' model"
$ws.Range("A17").Value = "44 u=4 6u6nPBC"
$ws.Range("A20").Value = "44 u=4 7u7nPBC"
$ws.Range("A23").Value = "44 u=4 8u8nPBC"
$ws.Range("H16").Value = "ED"
$ws.Range("B16").Value = "t'"
$ws.Range("A16").Value = "t' model"

# Row 16 - remaining headers
$ws.Range("D16").Value = "S.C. CPMC"
$ws.Range("E16").Value = "Err"
$ws.Range("F16").Value = "CPMC"
$ws.Range("G16").Value = "Err"

# Row 17-18 - "44 u=4 6u6nPBC"
$ws.Range("B17").Value = 0.2
$ws.Range("D17").Value = -19.223496000000001
$ws.Range("E17").Value = 0.0076816200000000001
$ws.Range("F17").Value = -19.189156195781042
$ws.Range("G17").Value = 0.005320384503046051
$ws.Range("H17").Value = -19.181403599999999

$ws.Range("B18").Value = -0.2
$ws.Range("D18").Value = -17.716363510000001
$ws.Range("E18").Value = 0.00499871
$ws.Range("F18").Value = -17.624560958986869
$ws.Range("G18").Value = 0.0055370825400400435
$ws.Range("H18").Value = -17.761400026

# Row 20-21 - "44 u=4 7u7nPBC"
$ws.Range("B20").Value = 0.2
$ws.Range("D20").Value = -17.092968259999999
$ws.Range("E20").Value = 0.0071702099999999998
$ws.Range("F20").Value = -16.773421572037179
$ws.Range("G20").Value = 0.0081332819481685827
$ws.Range("H20").Value = -17.182938

$ws.Range("B21").Value = -0.2
$ws.Range("D21").Value = -16.414027180000001
$ws.Range("E21").Value = 0.0056526299999999996
$ws.Range("F21").Value = -16.400135748483894
$ws.Range("G21").Value = 0.0058861681275161892
$ws.Range("H21").Value = -16.410562599999999

# Row 23-24 - "44 u=4 8u8nPBC"
$ws.Range("B23").Value = 0.2
$ws.Range("D23").Value = -13.503651550000001
$ws.Range("E23").Value = 0.01403358
$ws.Range("F23").Value = -13.237823021738489
$ws.Range("H23").Value = -13.627869

$ws.Range("B24").Value = -0.2
$ws.Range("D24").Value = -13.49667562
$ws.Range("E24").Value = 0.0089199699999999993
$ws.Range("F24").Value = -13.245576946580792
$ws.Range("G24").Value = 0.0082268898269285476
$ws.Range("H24").Value = -13.627869

# Set the formula after its dependent range (F23:F43, which includes F24) is populated
$ws.Range("G23").Formula = "=STDEV(F23:F43)/SQRT(20)"

# Selection state to match recorded view (scroll position isn't exposed via this COM surface)
$ws.Range("I28").Select()
